{"js": "// The document body contains a single 20-row x 5-column table whose cells\n// each hold one short arithmetic expression (e.g. \"23-5=\"). This edit swaps\n// every expression for a new one (same left-to-right, top-to-bottom cell\n// order), leaving the surrounding formatting (fonts, size, alignment, the\n// date heading, table structure, etc.) untouched.\nconst newValues = [\n  [\"90-55=\", \"48+49=\", \"9+67=\", \"33+49=\", \"9+59=\"],\n  [\"18+59=\", \"9+29=\", \"73-24=\", \"28+64=\", \"15+37=\"],\n  [\"84-57=\", \"58+17=\", \"29+65=\", \"59+23=\", \"48+26=\"],\n  [\"41-13=\", \"34+7=\", \"38+54=\", \"25-7=\", \"60-49=\"],\n  [\"57+19=\", \"83-26=\", \"57+4=\", \"54-18=\", \"17+37=\"],\n  [\"95-48=\", \"53+8=\", \"31-26=\", \"23-14=\", \"7+9=\"],\n  [\"92-5=\", \"58-19=\", \"26+69=\", \"74-6=\", \"15+18=\"],\n  [\"37-29=\", \"25+8=\", \"49+28=\", \"38+55=\", \"43+39=\"],\n  [\"42-14=\", \"25+39=\", \"53-26=\", \"33+29=\", \"53-5=\"],\n  [\"82-75=\", \"9+66=\", \"95-67=\", \"28+23=\", \"7+74=\"],\n  [\"66-39=\", \"55-49=\", \"55-37=\", \"71-17=\", \"72-54=\"],\n  [\"28+69=\", \"71-6=\", \"60-38=\", \"94-57=\", \"20-4=\"],\n  [\"83-24=\", \"94-16=\", \"73-17=\", \"94-6=\", \"32-5=\"],\n  [\"43-34=\", \"57+4=\", \"7+19=\", \"68+28=\", \"75-19=\"],\n  [\"43+38=\", \"26+57=\", \"43+28=\", \"9+72=\", \"35+29=\"],\n  [\"77+15=\", \"83-55=\", \"38+17=\", \"73+18=\", \"8+77=\"],\n  [\"61-53=\", \"70-54=\", \"6+78=\", \"19+63=\", \"35-28=\"],\n  [\"6+89=\", \"67-9=\", \"57-39=\", \"71-12=\", \"47+25=\"],\n  [\"57+37=\", \"40-3=\", \"38+26=\", \"61-29=\", \"55-39=\"],\n  [\"9+32=\", \"23+58=\", \"40-28=\", \"70-29=\", \"45+8=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\n// Build the replacement grid, preserving any row/column we don't expect\n// (defensive: only overwrite cells we have a mapped replacement for).\nconst current = table.values;\nconst updated = current.map((row, r) =>\n  row.map((cell, c) =>\n    (newValues[r] && newValues[r][c] !== undefined) ? newValues[r][c] : cell\n  )\n);\n\ntable.values = updated;\nawait context.sync();\n", "ps1": "# Update the 100 arithmetic-expression cells in the single 20x5 table to\n# their new values, in row-major (top-left -> bottom-right) order, matching\n# the order the cells appear in the document.\n$newValues = @(\n    '90-55=',\n    '48+49=',\n    '9+67=',\n    '33+49=',\n    '9+59=',\n    '18+59=',\n    '9+29=',\n    '73-24=',\n    '28+64=',\n    '15+37=',\n    '84-57=',\n    '58+17=',\n    '29+65=',\n    '59+23=',\n    '48+26=',\n    '41-13=',\n    '34+7=',\n    '38+54=',\n    '25-7=',\n    '60-49=',\n    '57+19=',\n    '83-26=',\n    '57+4=',\n    '54-18=',\n    '17+37=',\n    '95-48=',\n    '53+8=',\n    '31-26=',\n    '23-14=',\n    '7+9=',\n    '92-5=',\n    '58-19=',\n    '26+69=',\n    '74-6=',\n    '15+18=',\n    '37-29=',\n    '25+8=',\n    '49+28=',\n    '38+55=',\n    '43+39=',\n    '42-14=',\n    '25+39=',\n    '53-26=',\n    '33+29=',\n    '53-5=',\n    '82-75=',\n    '9+66=',\n    '95-67=',\n    '28+23=',\n    '7+74=',\n    '66-39=',\n    '55-49=',\n    '55-37=',\n    '71-17=',\n    '72-54=',\n    '28+69=',\n    '71-6=',\n    '60-38=',\n    '94-57=',\n    '20-4=',\n    '83-24=',\n    '94-16=',\n    '73-17=',\n    '94-6=',\n    '32-5=',\n    '43-34=',\n    '57+4=',\n    '7+19=',\n    '68+28=',\n    '75-19=',\n    '43+38=',\n    '26+57=',\n    '43+28=',\n    '9+72=',\n    '35+29=',\n    '77+15=',\n    '83-55=',\n    '38+17=',\n    '73+18=',\n    '8+77=',\n    '61-53=',\n    '70-54=',\n    '6+78=',\n    '19+63=',\n    '35-28=',\n    '6+89=',\n    '67-9=',\n    '57-39=',\n    '71-12=',\n    '47+25=',\n    '57+37=',\n    '40-3=',\n    '38+26=',\n    '61-29=',\n    '55-39=',\n    '9+32=',\n    '23+58=',\n    '40-28=',\n    '70-29=',\n    '45+8='\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$i = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n\nWrite-Output (\"Updated \" + $i + \" cells\")\n"}
